$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "69.547.04"
$ws.Range("E2").Value = "  -1.60%  "
$ws.Range("D3").Value = "3.680.08"
$ws.Range("E3").Value = "  -2.42%  "
$ws.Range("E4").Value = "  -0.22%  "
Set-TextValue $ws.Range("D5") "615.01"
$ws.Range("E5").Value = "  +0.27%  "
Set-TextValue $ws.Range("D6") "179.46"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").Value = "3.675.36"
$ws.Range("E7").Value = "  -2.41%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  -2.69%  "
$ws.Range("E10").Value = "  -3.10%  "
$ws.Range("E11").Value = "  -1.91%  "
Set-TextValue $ws.Range("D12") "0.479"
$ws.Range("E12").Value = "  -4.41%  "
Set-TextValue $ws.Range("D13") "39.88"
$ws.Range("E13").Value = "  -1.78%  "
$ws.Range("E14").Value = "  -3.27%  "
$ws.Range("D15").Value = "4.295.33"
$ws.Range("E15").Value = "  -2.45%  "
$ws.Range("D16").Value = "3.683.42"
$ws.Range("E16").Value = "  -2.79%  "
$ws.Range("D17").Value = "69.527.77"
$ws.Range("E17").Value = "  -2.00%  "
$ws.Range("E18").Value = "  -1.81%  "
$ws.Range("E19").Value = "  -0.82%  "
Set-TextValue $ws.Range("D20") "16.32"
$ws.Range("E20").Value = "  -3.50%  "
Set-TextValue $ws.Range("D21") "498.15"
$ws.Range("E21").Value = "  -5.18%  "
Set-TextValue $ws.Range("D22") "9.13"
$ws.Range("E22").Value = "  -3.18%  "
Set-TextValue $ws.Range("D23") "0.715"
$ws.Range("E23").Value = "  -4.24%  "
$ws.Range("E24").Value = "  -0.79%  "
Set-TextValue $ws.Range("D25") "86.09"
$ws.Range("E25").Value = "  -2.84%  "
$ws.Range("E28").Value = "  +3.36%  "
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("E30").Value = "  -2.79%  "
Set-TextValue $ws.Range("D31") "2.89"
$ws.Range("E31").Value = "  -0.62%  "
$ws.Range("E32").Value = "  -0.73%  "
Set-TextValue $ws.Range("D33") "30.06"
$ws.Range("E33").Value = "  -6.63%  "
$ws.Range("E34").Value = "  -1.75%  "
Set-TextValue $ws.Range("D35") "0.999"
$ws.Range("E35").Value = "  -0.40%  "
Set-TextValue $ws.Range("D36") "1.03"
$ws.Range("E36").Value = "  -1.18%  "
Set-TextValue $ws.Range("D37") "6.03"
$ws.Range("E37").Value = "  -2.29%  "
$ws.Range("E38").Value = "  +3.84%  "
Set-TextValue $ws.Range("D39") "0.339"
$ws.Range("E39").Value = "  -1.35%  "
Set-TextValue $ws.Range("D40") "49.92"
$ws.Range("E40").Value = "  -3.12%  "
$ws.Range("E41").Value = "  -6.94%  "
Set-TextValue $ws.Range("D42") "2.93"
$ws.Range("E42").Value = "  +4.50%  "
Set-TextValue $ws.Range("D43") "429.90"
$ws.Range("E43").Value = "  +1.14%  "
Set-TextValue $ws.Range("D44") "43.73"
$ws.Range("E44").Value = "  -0.81%  "
Set-TextValue $ws.Range("D45") "8.56"
$ws.Range("E45").Value = "  -3.56%  "
$ws.Range("D46").Value = "2.932.29"
$ws.Range("E46").Value = "  -7.05%  "
Set-TextValue $ws.Range("D47") "0.0358"
$ws.Range("E47").Value = "  -2.88%  "
Set-TextValue $ws.Range("D48") "27.40"
$ws.Range("E48").Value = "  -1.74%  "
$ws.Range("E49").Value = "  -0.05%  "
Set-TextValue $ws.Range("D50") "136.71"
$ws.Range("E50").Value = "  -3.21%  "
$ws.Range("E51").Value = "  -2.71%  "

# Rows 26 and 27 swap coin identity + values
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D26") "12.92"
$ws.Range("E26").Value = "  -4.93%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D27") "11.23"
$ws.Range("E27").Value = "  +2.59%  "

Write-Output "applied cryptos update"
